$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 459, shifting existing rows 459:481 down to 460:482
$ws.Rows.Item(459).Insert()

# Fill in the constant columns (same for every record in this block)
$ws.Cells.Item(459, 1).Value = 7
$ws.Cells.Item(459, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(459, 3).Value = "Ñuble"
$ws.Cells.Item(459, 4).Value = 45267
$ws.Cells.Item(459, 5).Value = 16
$ws.Cells.Item(459, 6).Value = "Fruta"
$ws.Cells.Item(459, 7).Value = 100103
$ws.Cells.Item(459, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(459, 9).Value = 100103006
$ws.Cells.Item(459, 10).Value = "Nectarín"
$ws.Cells.Item(459, 11).Value = "Big John"
$ws.Cells.Item(459, 12).Value = "Primera"
$ws.Cells.Item(459, 13).Value = 150
$ws.Cells.Item(459, 14).Value = 13000
$ws.Cells.Item(459, 15).Value = 13000
$ws.Cells.Item(459, 16).Value = 13000
$ws.Cells.Item(459, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(459, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(459, 19).Value = 867
$ws.Cells.Item(459, 20).Value = 15

# Match the date-cell number format used by the rest of column D
$ws.Cells.Item(459, 4).NumberFormat = $ws.Cells.Item(460, 4).NumberFormat
